# Update scripts wuth new tpm
# Refresh NATMI LR-pair TPM-derived metrics (Wnt5a-Fzd6) on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (FAPs -> ECs) ---
$ws.Range("I2").Value = 0.977669497583861
$ws.Range("J2").Value = 0.977669497583861
$ws.Range("M2").Value = 10.045207
$ws.Range("N2").Value = 30.135621
$ws.Range("O2").Value = 0.9365108453707793
$ws.Range("P2").Value = 0.9365108453707794
$ws.Range("Q2").Value = 75.48767133756499
$ws.Range("R2").Value = 679.389042038085
$ws.Range("S2").Value = 0.9155980876754868
$ws.Range("T2").Value = 0.9155980876754869

# --- Row 3 (FAPs -> FAPs) ---
$ws.Range("I3").Value = 0.977669497583861
$ws.Range("J3").Value = 0.977669497583861
$ws.Range("O3").Value = 0.03971513502725754
$ws.Range("P3").Value = 0.03971513502725754
$ws.Range("S3").Value = 0.03882827610857408
$ws.Range("T3").Value = 0.03882827610857408

# --- Row 4 (FAPs -> MuSCs) ---
$ws.Range("I4").Value = 0.977669497583861
$ws.Range("J4").Value = 0.977669497583861
$ws.Range("M4").Value = 0.255005
$ws.Range("N4").Value = 0.765015
$ws.Range("O4").Value = 0.02377401960196297
$ws.Range("P4").Value = 0.02377401960196297
$ws.Range("Q4").Value = 1.916310298975
$ws.Range("R4").Value = 17.246792690775
$ws.Range("S4").Value = 0.0232431337998
$ws.Range("T4").Value = 0.0232431337998

# --- Row 5 (MuSCs -> ECs) ---
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.171642
$ws.Range("H5").Value = 0.514926
$ws.Range("I5").Value = 0.02233050241613897
$ws.Range("J5").Value = 0.02233050241613898
$ws.Range("M5").Value = 10.045207
$ws.Range("N5").Value = 30.135621
$ws.Range("O5").Value = 0.9365108453707793
$ws.Range("P5").Value = 0.9365108453707794
$ws.Range("Q5").Value = 1.724179419894
$ws.Range("R5").Value = 15.517614779046
$ws.Range("S5").Value = 0.02091275769529254
$ws.Range("T5").Value = 0.02091275769529255

# --- Row 6 (MuSCs -> FAPs) ---
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.171642
$ws.Range("H6").Value = 0.514926
$ws.Range("I6").Value = 0.02233050241613897
$ws.Range("J6").Value = 0.02233050241613898
$ws.Range("O6").Value = 0.03971513502725754
$ws.Range("P6").Value = 0.03971513502725754
$ws.Range("Q6").Value = 0.073118233292
$ws.Range("R6").Value = 0.658064099628
$ws.Range("S6").Value = 0.0008868589186834601
$ws.Range("T6").Value = 0.0008868589186834602

# --- Row 7 (MuSCs -> MuSCs) ---
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.171642
$ws.Range("H7").Value = 0.514926
$ws.Range("I7").Value = 0.02233050241613897
$ws.Range("J7").Value = 0.02233050241613898
$ws.Range("M7").Value = 0.255005
$ws.Range("N7").Value = 0.765015
$ws.Range("O7").Value = 0.02377401960196297
$ws.Range("P7").Value = 0.02377401960196297
$ws.Range("Q7").Value = 0.04376956820999999
$ws.Range("R7").Value = 0.39392611389
$ws.Range("S7").Value = 0.0005308858021629693
$ws.Range("T7").Value = 0.0005308858021629694
